$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (row 1) to English snake_case names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Capitalize Spanish connector words (de, del, la, las, los, el, y) in
#    the state (col A) and municipality (col B) name cells, e.g.
#    "Pabellón de Arteaga" -> "Pabellón De Arteaga"
function FixName($text) {
    $res = $text
    $res = $res -replace '\bde\b', 'De'
    $res = $res -replace '\bdel\b', 'Del'
    $res = $res -replace '\bla\b', 'La'
    $res = $res -replace '\blas\b', 'Las'
    $res = $res -replace '\blos\b', 'Los'
    $res = $res -replace '\bel\b', 'El'
    $res = $res -replace '\by\b', 'Y'
    return $res
}

for ($row = 2; $row -le 1022; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $valA = $cellA.Value2
    if ($valA -ne $null -and $valA -ne "") {
        $cellA.Value = FixName($valA)
    }
    $cellB = $ws.Cells.Item($row, 2)
    $valB = $cellB.Value2
    if ($valB -ne $null -and $valB -ne "") {
        $cellB.Value = FixName($valB)
    }
}

# 3) Grand-total label row 1023: "TOTAL" -> "Total"
$ws.Range("A1023").Value = "Total"

# 4) Drop the trailing metadata/footnote rows (1025-1029): sample size,
#    source, author and date notes below the data table.
$ws.Range("A1025:A1029").EntireRow.Delete() | Out-Null
